$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 111; existing rows 111-183 shift down to 112-184,
# and the sheet dimension grows from A1:R183 to A1:R184.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new price record.
$ws.Cells.Item(111, 1).Value = 5
$ws.Cells.Item(111, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(111, 3).Value = "Maule"
$ws.Cells.Item(111, 4).Value = 44603
$ws.Cells.Item(111, 5).Value = 7
$ws.Cells.Item(111, 6).Value = 100112021
$ws.Cells.Item(111, 7).Value = "Ají"
$ws.Cells.Item(111, 8).Value = "Americana (o)"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 100
$ws.Cells.Item(111, 11).Value = 17000
$ws.Cells.Item(111, 12).Value = 17000
$ws.Cells.Item(111, 13).Value = 17000
$ws.Cells.Item(111, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(111, 15).Value = "Región del Maule"
$ws.Cells.Item(111, 16).Value = 680
$ws.Cells.Item(111, 17).Value = 25
$ws.Cells.Item(111, 18).Value = "Hortaliza"
